$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Factors")

# --- Expand Table2 on the Factors sheet with 2 new groups of 3 "conversion" columns ---
$tbl = $ws1.ListObjects.Item("Table2")
$tbl.Unlist()

# Insert 3 blank columns before the old "Human Hepatic Clearance (Cl_h)" column (G)
$ws1.Columns("G:I").Insert()
# Insert 3 blank columns before the old "Rag Hepatic Clearance (Cl_h)4" column
# (old column K, now shifted to N after the first insert above)
$ws1.Columns("N:P").Insert()

# Headers for the new columns
$ws1.Cells.Item(1,7).Value = "amount per gram liver"
$ws1.Cells.Item(1,8).Value = "gram liver per mL liver"
$ws1.Cells.Item(1,9).Value = "L liver per kg body weight"
$ws1.Cells.Item(1,14).Value = "amount per gram liver2"
$ws1.Cells.Item(1,15).Value = "gram liver per mL liver3"
$ws1.Cells.Item(1,16).Value = "L liver per kg body weight4"

# Row 2 (Cl_int_hep / hepatocyte)
$ws1.Range("G2").Formula = "=Constants!`$D`$4"
$ws1.Range("H2").Formula = "=Constants!`$D`$6"
$ws1.Range("I2").Formula = "=Constants!`$D`$5"
$ws1.Range("N2").Formula = "=Constants!`$F`$4"
$ws1.Range("O2").Formula = "=Constants!`$F`$6"
$ws1.Range("P2").Formula = "=Constants!`$F`$5"

# Row 3 (Cl_int_mic / microsome)
$ws1.Range("G3").Formula = "=Constants!`$D`$8"
$ws1.Range("H3").Formula = "=Constants!`$D`$6"
$ws1.Range("I3").Formula = "=Constants!`$D`$5"
$ws1.Range("N3").Formula = "=Constants!`$F`$8"
$ws1.Range("O3").Formula = "=Constants!`$F`$6"
$ws1.Range("P3").Formula = "=Constants!`$F`$5"

# Row 4 (Cl_int_hep / hepatocyte)
$ws1.Range("G4").Formula = "=Constants!`$D`$4"
$ws1.Range("H4").Formula = "=Constants!`$D`$6"
$ws1.Range("I4").Formula = "=Constants!`$D`$5"
$ws1.Range("N4").Formula = "=Constants!`$F`$4"
$ws1.Range("O4").Formula = "=Constants!`$F`$6"
$ws1.Range("P4").Formula = "=Constants!`$F`$5"

# Row 5 (Cl_int_mic / microsome)
$ws1.Range("G5").Formula = "=Constants!`$D`$8"
$ws1.Range("H5").Formula = "=Constants!`$D`$6"
$ws1.Range("I5").Formula = "=Constants!`$D`$5"
$ws1.Range("N5").Formula = "=Constants!`$F`$8"
$ws1.Range("O5").Formula = "=Constants!`$F`$6"
$ws1.Range("P5").Formula = "=Constants!`$F`$5"

# Re-point the "Human/Rag Hepatic Clearance" formulas to the new per-row
# conversion columns instead of directly to the Constants sheet.
$ws1.Range("J2").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/10^6*60*1000"
$ws1.Range("Q2").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/10^6*60*1000"

$ws1.Range("J3").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/1000"
$ws1.Range("Q3").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/1000"

$ws1.Range("J4").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/10^3*1000"
$ws1.Range("Q4").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/10^3*1000"

$ws1.Range("J5").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/1000/1000*60"
$ws1.Range("Q5").Formula = "=Table2[[#This Row],[Value]]*Table2[[#This Row],[amount per gram liver]]*Table2[[#This Row],[gram liver per mL liver]]*Table2[[#This Row],[L liver per kg body weight]]/1000/1000*60"

# Rebuild the table over the new A1:R5 range so headers / columns realign.
$tbl2 = $ws1.ListObjects.Add(1, $ws1.Range("A1:R5"), 0, 1)
$tbl2.Name = "Table2"

# Apply the 0.0000 number format to the new "L liver per kg body weight4" column
$ws1.Range("P2:P5").NumberFormat = "0.0000"
